$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.050945333333334
$ws.Range("H2").Value = 15.152836
$ws.Range("I2").Value = 0.6229573853973969
$ws.Range("J2").Value = 0.6229573853973969
$ws.Range("M2").Value = 28.25464766666667
$ws.Range("N2").Value = 84.763943
$ws.Range("O2").Value = 0.3168758800036845
$ws.Range("P2").Value = 0.3168758800036845
$ws.Range("Q2").Value = 142.7126807769276
$ws.Range("R2").Value = 1284.414126992348
$ws.Range("S2").Value = 0.1974001697025946
$ws.Range("T2").Value = 0.1974001697025946

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.050945333333334
$ws.Range("H3").Value = 15.152836
$ws.Range("I3").Value = 0.6229573853973969
$ws.Range("J3").Value = 0.6229573853973969
$ws.Range("M3").Value = 0.04671833333333333
$ws.Range("N3").Value = 0.140155
$ws.Range("O3").Value = 0.0005239461189519747
$ws.Range("P3").Value = 0.0005239461189519747
$ws.Range("Q3").Value = 0.2359717477311111
$ws.Range("R3").Value = 2.12374572958
$ws.Range("S3").Value = 0.0003263961043514357
$ws.Range("T3").Value = 0.0003263961043514357

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.050945333333334
$ws.Range("H4").Value = 15.152836
$ws.Range("I4").Value = 0.6229573853973969
$ws.Range("J4").Value = 0.6229573853973969
$ws.Range("M4").Value = 60.86492733333333
$ws.Range("N4").Value = 182.594782
$ws.Range("O4").Value = 0.6826001738773636
$ws.Range("P4").Value = 0.6826001738773636
$ws.Range("Q4").Value = 307.4254206779725
$ws.Range("R4").Value = 2766.828786101752
$ws.Range("S4").Value = 0.4252308195904509
$ws.Range("T4").Value = 0.4252308195904509

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.696109
$ws.Range("H5").Value = 8.088327
$ws.Range("I5").Value = 0.3325240925302148
$ws.Range("J5").Value = 0.3325240925302149
$ws.Range("M5").Value = 28.25464766666667
$ws.Range("N5").Value = 84.763943
$ws.Range("O5").Value = 0.3168758800036845
$ws.Range("P5").Value = 0.3168758800036845
$ws.Range("Q5").Value = 76.177609865929
$ws.Range("R5").Value = 685.598488793361
$ws.Range("S5").Value = 0.1053688644429384
$ws.Range("T5").Value = 0.1053688644429384

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.696109
$ws.Range("H6").Value = 8.088327
$ws.Range("I6").Value = 0.3325240925302148
$ws.Range("J6").Value = 0.3325240925302149
$ws.Range("M6").Value = 0.04671833333333333
$ws.Range("N6").Value = 0.140155
$ws.Range("O6").Value = 0.0005239461189519747
$ws.Range("P6").Value = 0.0005239461189519747
$ws.Range("Q6").Value = 0.125957718965
$ws.Range("R6").Value = 1.133619470685
$ws.Range("S6").Value = 0.0001742247077392334
$ws.Range("T6").Value = 0.0001742247077392334

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.696109
$ws.Range("H7").Value = 8.088327
$ws.Range("I7").Value = 0.3325240925302148
$ws.Range("J7").Value = 0.3325240925302149
$ws.Range("M7").Value = 60.86492733333333
$ws.Range("N7").Value = 182.594782
$ws.Range("O7").Value = 0.6826001738773636
$ws.Range("P7").Value = 0.6826001738773636
$ws.Range("Q7").Value = 164.098478367746
$ws.Range("R7").Value = 1476.886305309714
$ws.Range("S7").Value = 0.2269810033795372
$ws.Range("T7").Value = 0.2269810033795372

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3609566666666666
$ws.Range("H8").Value = 1.08287
$ws.Range("I8").Value = 0.04451852207238823
$ws.Range("J8").Value = 0.04451852207238824
$ws.Range("M8").Value = 28.25464766666667
$ws.Range("N8").Value = 84.763943
$ws.Range("O8").Value = 0.3168758800036845
$ws.Range("P8").Value = 0.3168758800036845
$ws.Range("Q8").Value = 10.19870343960111
$ws.Range("R8").Value = 91.78833095640999
$ws.Range("S8").Value = 0.01410684585815147
$ws.Range("T8").Value = 0.01410684585815148

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3609566666666666
$ws.Range("H9").Value = 1.08287
$ws.Range("I9").Value = 0.04451852207238823
$ws.Range("J9").Value = 0.04451852207238824
$ws.Range("M9").Value = 0.04671833333333333
$ws.Range("N9").Value = 0.140155
$ws.Range("O9").Value = 0.0005239461189519747
$ws.Range("P9").Value = 0.0005239461189519747
$ws.Range("Q9").Value = 0.01686329387222222
$ws.Range("R9").Value = 0.15176964485
$ws.Range("S9").Value = 0.00002332530686130564
$ws.Range("T9").Value = 0.00002332530686130564

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3609566666666666
$ws.Range("H10").Value = 1.08287
$ws.Range("I10").Value = 0.04451852207238823
$ws.Range("J10").Value = 0.04451852207238824
$ws.Range("M10").Value = 60.86492733333333
$ws.Range("N10").Value = 182.594782
$ws.Range("O10").Value = 0.6826001738773636
$ws.Range("P10").Value = 0.6826001738773636
$ws.Range("Q10").Value = 21.96960128714889
$ws.Range("R10").Value = 197.72641158434
$ws.Range("S10").Value = 0.03038835090737546
$ws.Range("T10").Value = 0.03038835090737546
